## Add the 20 new "Linguistics" user-story rows (rows 42-61) to Foglio1,
## mirroring the commit "Add US related to Linguistics domain".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1) Seed every new shared string in the exact first-occurrence order the
#    original author used, so the appended <si> entries land at the same
#    indices (70 = "Linguistics" Domain, 71-90 = the 20 new user stories).
# ----------------------------------------------------------------------
$ws.Range("C42").Value = "Linguistics"

# Column E (user story) written in first-occurrence order - note rows 52/53
# are written out of sheet order (53 before 52) to reproduce the original
# shared-string sequence exactly.
$ws.Range("E42").Value = "As a linguist, I want to employ adversarial learning techniques to detect and mitigate dialectal biases in automated transcription systems, ensuring accurate representation of regional language variations."
$ws.Range("E43").Value = "As a linguist, I want to utilize CNNs to automatically analyze and categorize handwritten manuscripts based on historical writing styles and script variations, aiding in the classification and digitization of ancient texts."
$ws.Range("E44").Value = "As a linguist, I want to create a conversational agent capable of simulating dialogues in multiple languages, allowing users to practice conversational fluency and cultural expressions in a natural setting."
$ws.Range("E45").Value = "As a linguist, I want to apply decision tree algorithms to categorize language samples into phonetic groups, facilitating the study of phonological evolution and dialectal variations over time."
$ws.Range("E46").Value = "As a linguist, I want to apply document classification techniques to categorize language corpora into historical periods, facilitating chronological studies of language evolution and linguistic changes."
$ws.Range("E47").Value = "As a linguist, I want to develop entity extraction algorithms to automatically detect and categorize named entities such as proper names and geographic locations in historical texts, aiding in geographical and historical linguistic research."
$ws.Range("E48").Value = "As a linguist, I want to use feature selection techniques to optimize the selection of discourse markers and pragmatic features in spoken dialogues, enhancing the analysis of conversational strategies and communicative intentions."
$ws.Range("E49").Value = "As a linguist, I want to address imbalanced datasets in historical language corpora by applying techniques that balance the representation of under-documented languages, ensuring fair and comprehensive linguistic analysis."
$ws.Range("E50").Value = "As a linguist, I want to employ keyword extraction techniques on multilingual dictionaries to identify semantic clusters and cross-linguistic equivalences, supporting comparative lexical studies."
$ws.Range("E51").Value = "As a linguist, I want to utilize k-NN models to identify similar linguistic patterns across different languages based on shared morphological and syntactic features, supporting typological studies."
$ws.Range("E53").Value = "As a linguist, I want to use neural networks for automatic language identification in multilingual texts, improving the efficiency of language documentation and corpus compilation efforts."
$ws.Range("E52").Value = "As a linguist, I want to use multi-label classification algorithms to categorize language learning materials into proficiency levels (beginner, intermediate, advanced) and language skills (listening, speaking, reading, writing)."
$ws.Range("E54").Value = "As a linguist, I want to use a random forest algorithm to classify texts based on syntactic structures, enabling automated categorization of sentences into grammatical types (e.g., declarative, interrogative, imperative)."
$ws.Range("E55").Value = "As a linguist, I want to use semantic similarity algorithms to identify synonymous terms in multilingual dictionaries, facilitating the compilation of comprehensive lexical resources."
$ws.Range("E56").Value = "As a linguist, I want to use sentiment analysis techniques to analyze student feedback on language courses and teaching methodologies, improving instructional practices and student satisfaction."
$ws.Range("E57").Value = "As a linguist, I want to use a speech to text system to transcribe language learning sessions, allowing for detailed analysis of pronunciation and speech patterns to enhance teaching methods."
$ws.Range("E58").Value = "As a linguist, I want to use text categorization algorithms to classify research papers into linguistic subfields such as phonetics, syntax, and semantics, facilitating targeted literature reviews."
$ws.Range("E59").Value = "As a linguist, I want to apply unsupervised clustering algorithms to categorize phonetic data from different dialects, identifying distinct phonological patterns and variations."
$ws.Range("E60").Value = "As a linguist, I want to use voice recognition technology to create language learning tools that provide real-time feedback on pronunciation accuracy and intonation."
$ws.Range("E61").Value = "As a linguist, I want to use word embedding techniques to map words from different languages into a shared semantic space, allowing for cross-linguistic comparison of lexical semantics and conceptual structures."

# ----------------------------------------------------------------------
# 2) Fill in the remaining columns for every new row (A = domain cluster,
#    B = constant weight, C = domain, D = ML task, F = prompt label). All
#    of these reuse shared strings that already exist, so order is free.
# ----------------------------------------------------------------------
$ws.Range("A42").Value = "Literature & Linguistics"
$ws.Range("B42").Value = 5
$ws.Range("C42").Value = "Linguistics"
$ws.Range("D42").Value = "adversarial learning"
$ws.Range("F42").Value = "Domain_FSPrompt"
$ws.Range("A43").Value = "Literature & Linguistics"
$ws.Range("B43").Value = 5
$ws.Range("C43").Value = "Linguistics"
$ws.Range("D43").Value = "cnn"
$ws.Range("F43").Value = "Domain_FSPrompt"
$ws.Range("A44").Value = "Literature & Linguistics"
$ws.Range("B44").Value = 5
$ws.Range("C44").Value = "Linguistics"
$ws.Range("D44").Value = "conversational agent"
$ws.Range("F44").Value = "Domain_FSPrompt"
$ws.Range("A45").Value = "Literature & Linguistics"
$ws.Range("B45").Value = 5
$ws.Range("C45").Value = "Linguistics"
$ws.Range("D45").Value = "decision tree"
$ws.Range("F45").Value = "Domain_FSPrompt"
$ws.Range("A46").Value = "Literature & Linguistics"
$ws.Range("B46").Value = 5
$ws.Range("C46").Value = "Linguistics"
$ws.Range("D46").Value = "document classification"
$ws.Range("F46").Value = "Domain_FSPrompt"
$ws.Range("A47").Value = "Literature & Linguistics"
$ws.Range("B47").Value = 5
$ws.Range("C47").Value = "Linguistics"
$ws.Range("D47").Value = "entity extraction"
$ws.Range("F47").Value = "Domain_FSPrompt"
$ws.Range("A48").Value = "Literature & Linguistics"
$ws.Range("B48").Value = 5
$ws.Range("C48").Value = "Linguistics"
$ws.Range("D48").Value = "feature selection"
$ws.Range("F48").Value = "Domain_FSPrompt"
$ws.Range("A49").Value = "Literature & Linguistics"
$ws.Range("B49").Value = 5
$ws.Range("C49").Value = "Linguistics"
$ws.Range("D49").Value = "imbalanced dataset"
$ws.Range("F49").Value = "Domain_FSPrompt"
$ws.Range("A50").Value = "Literature & Linguistics"
$ws.Range("B50").Value = 5
$ws.Range("C50").Value = "Linguistics"
$ws.Range("D50").Value = "keyword extraction"
$ws.Range("F50").Value = "Domain_FSPrompt"
$ws.Range("A51").Value = "Literature & Linguistics"
$ws.Range("B51").Value = 5
$ws.Range("C51").Value = "Linguistics"
$ws.Range("D51").Value = "k-nearest neighbor"
$ws.Range("F51").Value = "Domain_FSPrompt"
$ws.Range("A52").Value = "Literature & Linguistics"
$ws.Range("B52").Value = 5
$ws.Range("C52").Value = "Linguistics"
$ws.Range("D52").Value = "multi-label classification"
$ws.Range("F52").Value = "Domain_FSPrompt"
$ws.Range("A53").Value = "Literature & Linguistics"
$ws.Range("B53").Value = 5
$ws.Range("C53").Value = "Linguistics"
$ws.Range("D53").Value = "neural network"
$ws.Range("F53").Value = "Domain_FSPrompt"
$ws.Range("A54").Value = "Literature & Linguistics"
$ws.Range("B54").Value = 5
$ws.Range("C54").Value = "Linguistics"
$ws.Range("D54").Value = "random forest"
$ws.Range("F54").Value = "Domain_FSPrompt"
$ws.Range("A55").Value = "Literature & Linguistics"
$ws.Range("B55").Value = 5
$ws.Range("C55").Value = "Linguistics"
$ws.Range("D55").Value = "semantic similarity"
$ws.Range("F55").Value = "Domain_FSPrompt"
$ws.Range("A56").Value = "Literature & Linguistics"
$ws.Range("B56").Value = 5
$ws.Range("C56").Value = "Linguistics"
$ws.Range("D56").Value = "sentiment analysis"
$ws.Range("F56").Value = "Domain_FSPrompt"
$ws.Range("A57").Value = "Literature & Linguistics"
$ws.Range("B57").Value = 5
$ws.Range("C57").Value = "Linguistics"
$ws.Range("D57").Value = "speech to text"
$ws.Range("F57").Value = "Domain_FSPrompt"
$ws.Range("A58").Value = "Literature & Linguistics"
$ws.Range("B58").Value = 5
$ws.Range("C58").Value = "Linguistics"
$ws.Range("D58").Value = "text categorization"
$ws.Range("F58").Value = "Domain_FSPrompt"
$ws.Range("A59").Value = "Literature & Linguistics"
$ws.Range("B59").Value = 5
$ws.Range("C59").Value = "Linguistics"
$ws.Range("D59").Value = "unsupervised clustering"
$ws.Range("F59").Value = "Domain_FSPrompt"
$ws.Range("A60").Value = "Literature & Linguistics"
$ws.Range("B60").Value = 5
$ws.Range("C60").Value = "Linguistics"
$ws.Range("D60").Value = "voice recognition"
$ws.Range("F60").Value = "Domain_FSPrompt"
$ws.Range("A61").Value = "Literature & Linguistics"
$ws.Range("B61").Value = 5
$ws.Range("C61").Value = "Linguistics"
$ws.Range("D61").Value = "word embedding"
$ws.Range("F61").Value = "Domain_FSPrompt"

# ----------------------------------------------------------------------
# 3) Match the per-column formatting used throughout the sheet (bold/blue
#    A:C, bold D, plain E:F) by copying an existing data row's formats -
#    the E column is intentionally left at its column default format.
# ----------------------------------------------------------------------
$ws.Range("A2:D2").Copy()
$ws.Range("A42:D42").PasteSpecial(-4122)
$ws.Range("A43:D43").PasteSpecial(-4122)
$ws.Range("A44:D44").PasteSpecial(-4122)
$ws.Range("A45:D45").PasteSpecial(-4122)
$ws.Range("A46:D46").PasteSpecial(-4122)
$ws.Range("A47:D47").PasteSpecial(-4122)
$ws.Range("A48:D48").PasteSpecial(-4122)
$ws.Range("A49:D49").PasteSpecial(-4122)
$ws.Range("A50:D50").PasteSpecial(-4122)
$ws.Range("A51:D51").PasteSpecial(-4122)
$ws.Range("A52:D52").PasteSpecial(-4122)
$ws.Range("A53:D53").PasteSpecial(-4122)
$ws.Range("A54:D54").PasteSpecial(-4122)
$ws.Range("A55:D55").PasteSpecial(-4122)
$ws.Range("A56:D56").PasteSpecial(-4122)
$ws.Range("A57:D57").PasteSpecial(-4122)
$ws.Range("A58:D58").PasteSpecial(-4122)
$ws.Range("A59:D59").PasteSpecial(-4122)
$ws.Range("A60:D60").PasteSpecial(-4122)
$ws.Range("A61:D61").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F42").PasteSpecial(-4122)
$ws.Range("F43").PasteSpecial(-4122)
$ws.Range("F44").PasteSpecial(-4122)
$ws.Range("F45").PasteSpecial(-4122)
$ws.Range("F46").PasteSpecial(-4122)
$ws.Range("F47").PasteSpecial(-4122)
$ws.Range("F48").PasteSpecial(-4122)
$ws.Range("F49").PasteSpecial(-4122)
$ws.Range("F50").PasteSpecial(-4122)
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("F52").PasteSpecial(-4122)
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("F54").PasteSpecial(-4122)
$ws.Range("F55").PasteSpecial(-4122)
$ws.Range("F56").PasteSpecial(-4122)
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("F58").PasteSpecial(-4122)
$ws.Range("F59").PasteSpecial(-4122)
$ws.Range("F60").PasteSpecial(-4122)
$ws.Range("F61").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# 4) Scroll/select so the view lands on the newly-added data, matching the
#    saved sheet view (topLeftCell A45, active cell E60).
# ----------------------------------------------------------------------
$excel.Goto($ws.Range("A45"), $true)
$ws.Range("E60").Select()

